# 🚌 141: 30/12 16:37 LP1912+6203+6173
# Append newly-scraped rows to the three schedule sheets and refresh the
# "Última actualización" / "Total filas" banner cells.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "LP1912": columns A(blank) B=Hora_Scrap C=Hora_Llegada D=Línea
#                 E=Minutos(number) F=Parada G=Fecha
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 30/12/2025 13:37:49"
$ws1.Range("A3").Value = "Total filas: 260"

$lp1912Rows = @(
    @("13:37:39", "13:46", "16_SANTA ANA", 9, "LP1912", "30/12/2025"),
    @("13:37:39", "13:56", "16_P MOR-167 Y 521", 19, "LP1912", "30/12/2025"),
    @("13:37:39", "14:04", "17_ROMERO", 27, "LP1912", "30/12/2025"),
    @("13:37:39", "14:04", "23_HERNANDEZ", 27, "LP1912", "30/12/2025"),
    @("13:37:39", "14:06", "16_SANTA ANA", 29, "LP1912", "30/12/2025"),
    @("13:37:39", "14:21", "26_HERNANDEZ", 44, "LP1912", "30/12/2025"),
    @("13:37:39", "14:44", "14_ABASTO", 67, "LP1912", "30/12/2025"),
    @("13:37:39", "14:56", "16_P MOR-SANTA ANA", 79, "LP1912", "30/12/2025"),
    @("13:37:39", "14:58", "215B_EL PATO", 81, "LP1912", "30/12/2025"),
    @("13:37:39", "15:00", "81_EL PELIGRO", 83, "LP1912", "30/12/2025"),
    @("13:37:39", "15:05", "10_OLMOS", 88, "LP1912", "30/12/2025")
)

$startRow = 251
for ($i = 0; $i -lt $lp1912Rows.Count; $i++) {
    $r = $startRow + $i
    $row = $lp1912Rows[$i]
    $ws1.Cells.Item($r, 2).Value = $row[0]
    $ws1.Cells.Item($r, 3).Value = $row[1]
    $ws1.Cells.Item($r, 4).Value = $row[2]
    $ws1.Cells.Item($r, 5).Value = $row[3]
    $ws1.Cells.Item($r, 6).Value = $row[4]
    $ws1.Cells.Item($r, 7).Value = $row[5]
}

# ---------------------------------------------------------------------
# Sheet "LP1912-215": columns A(blank) B=Fecha C=Hora_Scrap D=Hora_Llegada
#                      E=Línea F=Minutos(number) G=Parada
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 30/12/2025 13:37:49"
$ws2.Range("A3").Value = "Total filas: 19"

$ws2.Cells.Item(20, 2).Value = "30/12/2025"
$ws2.Cells.Item(20, 3).Value = "13:37:39"
$ws2.Cells.Item(20, 4).Value = "14:58"
$ws2.Cells.Item(20, 5).Value = "215B_EL PATO"
$ws2.Cells.Item(20, 6).Value = 81
$ws2.Cells.Item(20, 7).Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet "6203-6173": columns A(blank) B=Fecha C=Hora_Scrap D=Hora_Llegada
#                     E=Línea F=Minutos(number) G=Parada
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 30/12/2025 13:37:49"
$ws3.Range("A3").Value = "Total filas: 37"

$ws3.Cells.Item(37, 2).Value = "30/12/2025"
$ws3.Cells.Item(37, 3).Value = "13:37:49"
$ws3.Cells.Item(37, 4).Value = "14:09"
$ws3.Cells.Item(37, 5).Value = "215A_LA PLATA"
$ws3.Cells.Item(37, 6).Value = 32
$ws3.Cells.Item(37, 7).Value = "L6173"

$ws3.Cells.Item(38, 2).Value = "30/12/2025"
$ws3.Cells.Item(38, 3).Value = "13:37:44"
$ws3.Cells.Item(38, 4).Value = "14:52"
$ws3.Cells.Item(38, 5).Value = "215D_LA PLATA"
$ws3.Cells.Item(38, 6).Value = 75
$ws3.Cells.Item(38, 7).Value = "L6203"
